$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 157, shifting existing rows 157-223 down to 160-226
$ws.Rows("157:159").Insert()

# Row 157
$ws.Range("A157").Value = 3
$ws.Range("B157").Value = 'Femacal de La Calera'
$ws.Range("C157").Value = 'Coquimbo'
$ws.Range("D157").Value = 44917
$ws.Range("E157").Value = 5
$ws.Range("F157").Value = 'Fruta'
$ws.Range("G157").Value = 100103
$ws.Range("H157").Value = 'Frutos de hueso (carozo)'
$ws.Range("I157").Value = 100103002
$ws.Range("J157").Value = 'Ciruela'
$ws.Range("K157").Value = 'Black Amber'
$ws.Range("L157").Value = 'Especial'
$ws.Range("M157").Value = 75
$ws.Range("N157").Value = 17000
$ws.Range("O157").Value = 17000
$ws.Range("P157").Value = 17000
$ws.Range("Q157").Value = '$/caja 15 kilos granel'
$ws.Range("R157").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S157").Value = 1133
$ws.Range("T157").Value = 15

# Row 158
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = 'Femacal de La Calera'
$ws.Range("C158").Value = 'Coquimbo'
$ws.Range("D158").Value = 44917
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = 'Fruta'
$ws.Range("G158").Value = 100103
$ws.Range("H158").Value = 'Frutos de hueso (carozo)'
$ws.Range("I158").Value = 100103002
$ws.Range("J158").Value = 'Ciruela'
$ws.Range("K158").Value = 'Black Amber'
$ws.Range("L158").Value = 'Primera'
$ws.Range("M158").Value = 85
$ws.Range("N158").Value = 14000
$ws.Range("O158").Value = 14000
$ws.Range("P158").Value = 14000
$ws.Range("Q158").Value = '$/caja 15 kilos granel'
$ws.Range("R158").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S158").Value = 933
$ws.Range("T158").Value = 15

# Row 159
$ws.Range("A159").Value = 3
$ws.Range("B159").Value = 'Femacal de La Calera'
$ws.Range("C159").Value = 'Coquimbo'
$ws.Range("D159").Value = 44917
$ws.Range("E159").Value = 5
$ws.Range("F159").Value = 'Fruta'
$ws.Range("G159").Value = 100103
$ws.Range("H159").Value = 'Frutos de hueso (carozo)'
$ws.Range("I159").Value = 100103002
$ws.Range("J159").Value = 'Ciruela'
$ws.Range("K159").Value = 'Black Amber'
$ws.Range("L159").Value = 'Segunda'
$ws.Range("M159").Value = 80
$ws.Range("N159").Value = 12000
$ws.Range("O159").Value = 12000
$ws.Range("P159").Value = 12000
$ws.Range("Q159").Value = '$/caja 15 kilos granel'
$ws.Range("R159").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S159").Value = 800
$ws.Range("T159").Value = 15

